$d = $word.ActiveDocument

# --- Core content fix: correct the default value of OMSAssistEnable from TRUE to FALSE ---
$t = $d.Tables.Item(2)
$cell = $t.Cell(10, 4)
$r = $cell.Range
# Trim the trailing end-of-cell marker so we only touch the visible text.
$r.End = $r.End - 1
$r.Text = "FALSE"

# --- Move the "_GoBack" bookmark so it marks this last-edited spot (matches Word's ---
# --- behaviour of recording the most recent edit location) instead of the document end. ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$cell2 = $t.Cell(10, 4)
$r2 = $cell2.Range
$r2.End = $r2.End - 1
$d.Bookmarks.Add("_GoBack", $r2)

Write-Host "Edit complete"
